$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Fund name -> "SAAS Fund" for all data rows (2-7) ---
$ws.Range("A2:A7").Value = "SAAS Fund"

# --- Column B: Investor names -> Investor 1..4 pattern ---
$ws.Cells.Item(2,2).Value = "Investor 1"
$ws.Cells.Item(3,2).Value = "Investor 2"
$ws.Cells.Item(4,2).Value = "Investor 3"
$ws.Cells.Item(5,2).Value = "Investor 4"
$ws.Cells.Item(6,2).Value = "Investor 1"
$ws.Cells.Item(7,2).Value = "Investor 2"

# --- Column C: Capital Distribution label -> Distribution 1 / Distribution 2 (with variant spacing) ---
$ws.Cells.Item(2,3).Value = "Distribution 1 "
$ws.Cells.Item(3,3).Value = "Distribution 1     "
$ws.Cells.Item(4,3).Value = "Distribution 1"
$ws.Cells.Item(5,3).Value = "Distribution 1"
$ws.Cells.Item(6,3).Value = "Distribution 2"
$ws.Cells.Item(7,3).Value = "Distribution 2"

# --- Column H: new "Folio No" data, rows 2-6 use a freshly introduced style (Arial 11) ---
$c = $ws.Cells.Item(2,8)
$c.Value = 6
$c.Font.Name = "Arial"
$c.Font.Size = 11

$c = $ws.Cells.Item(3,8)
$c.Value = 7
$c.Font.Name = "Arial"
$c.Font.Size = 11

$c = $ws.Cells.Item(4,8)
$c.Value = 8
$c.Font.Name = "Arial"
$c.Font.Size = 11

$c = $ws.Cells.Item(5,8)
$c.Value = 9
$c.Font.Name = "Arial"
$c.Font.Size = 11

$c = $ws.Cells.Item(6,8)
$c.Value = 6
$c.Font.Name = "Arial"
$c.Font.Size = 11

# Row 7 keeps the plain pre-existing data style (no explicit font change)
$ws.Cells.Item(7,8).Value = 7

# --- Update the saved selection / active cell ---
$ws.Range("C4").Select() | Out-Null
